# Populate the manufacturer / year / reference_number / model_name columns
# (C, D, E, F) for the lot rows that were scraped, and touch the other
# per-lot detail columns (L, M, O, P, Q, R, S, T) that the scraper always
# writes (even when empty) so every lot row has a full set of columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlLineStyleNone - used as a harmless "touch" so a blank cell is written
# out to the sheet (a plain empty Value assignment clears/removes the
# cell instead of leaving it present-but-blank).
$xlLineStyleNone = -4142

function Set-BlankCell([string]$addr) {
    $ws.Range($addr).Borders.LineStyle = $xlLineStyleNone
}

# --- Row 2: NATIONAL WATCH ---------------------------------------------
$ws.Range("C2").Value = "NATIONAL WATCH"
$ws.Range("F2").Value = "NATIONAL WATCH"
Set-BlankCell "D2"
Set-BlankCell "E2"
Set-BlankCell "L2"
Set-BlankCell "M2"
Set-BlankCell "O2"
Set-BlankCell "P2"
Set-BlankCell "Q2"
Set-BlankCell "R2"
Set-BlankCell "S2"
Set-BlankCell "T2"

# --- Row 3: GALLET -------------------------------------------------------
$ws.Range("C3").Value = "GALLET"
$ws.Range("F3").Value = "GALLET"
Set-BlankCell "D3"
Set-BlankCell "E3"
Set-BlankCell "L3"
Set-BlankCell "M3"
Set-BlankCell "O3"
Set-BlankCell "P3"
Set-BlankCell "Q3"
Set-BlankCell "R3"
Set-BlankCell "S3"
Set-BlankCell "T3"

# --- Row 4: HELVETIA -------------------------------------------------------
$ws.Range("C4").Value = "HELVETIA"
$ws.Range("F4").Value = "HELVETIA"
Set-BlankCell "D4"
Set-BlankCell "E4"
Set-BlankCell "L4"
Set-BlankCell "M4"
Set-BlankCell "O4"
Set-BlankCell "P4"
Set-BlankCell "Q4"
Set-BlankCell "R4"
Set-BlankCell "S4"
Set-BlankCell "T4"

# --- Row 5: BELL & ROSS ----------------------------------------------------
$ws.Range("C5").Value = "BELL & ROSS REF. BR 01-97 PVD STEEL LIMITED EDITION`nBell & Ross"
$ws.Range("E5").Value = "BR"
$ws.Range("F5").Value = "BELL & ROSS REF. BR 01-97 PVD STEEL LIMITED EDITION`nBell & Ross"
Set-BlankCell "D5"
Set-BlankCell "L5"
Set-BlankCell "M5"
Set-BlankCell "O5"
Set-BlankCell "P5"
Set-BlankCell "Q5"
Set-BlankCell "R5"
Set-BlankCell "S5"
Set-BlankCell "T5"

# --- Row 6: CHRONOSWISS ------------------------------------------------
$ws.Range("C6").Value = "CHRONOSWISS KLASSIK REF. CH 7443 CHRONOGRAPH STEEL `nChronoswiss"
$ws.Range("D6").Value = "7443"
$ws.Range("E6").Value = "CH"
$ws.Range("F6").Value = "CHRONOSWISS KLASSIK REF. CH 7443 CHRONOGRAPH STEEL `nChronoswiss"
Set-BlankCell "L6"
Set-BlankCell "M6"
Set-BlankCell "O6"
Set-BlankCell "P6"
Set-BlankCell "Q6"
Set-BlankCell "R6"
Set-BlankCell "S6"
Set-BlankCell "T6"

# --- Row 7: VACHERON & CONSTANTIN --------------------------------------
$ws.Range("C7").Value = "VACHERON & CONSTANTIN REF. 33093 YELLOW GOLD`nVacheron & Constantin"
$ws.Range("E7").Value = "33093"
$ws.Range("F7").Value = "VACHERON & CONSTANTIN REF. 33093 YELLOW GOLD`nVacheron & Constantin"
Set-BlankCell "D7"
Set-BlankCell "L7"
Set-BlankCell "M7"
Set-BlankCell "O7"
Set-BlankCell "P7"
Set-BlankCell "Q7"
Set-BlankCell "R7"
Set-BlankCell "S7"
Set-BlankCell "T7"

# --- Row 8: ORBITA 6 ROTOR ----------------------------------------------
$ws.Range("C8").Value = "ORBITA 6 ROTOR WATCH `nWINDING CABINET WOOD `nOrbita"
$ws.Range("F8").Value = "ORBITA 6 ROTOR WATCH `nWINDING CABINET WOOD `nOrbita"
Set-BlankCell "D8"
Set-BlankCell "E8"
Set-BlankCell "L8"
Set-BlankCell "M8"
Set-BlankCell "O8"
Set-BlankCell "P8"
Set-BlankCell "Q8"
Set-BlankCell "R8"
Set-BlankCell "S8"
Set-BlankCell "T8"

# --- Row 9: ORBITA 12 ROTOR ---------------------------------------------
$ws.Range("C9").Value = "ORBITA 12 ROTOR WATCH `nWINDING CABINET WOOD `nOrbita"
$ws.Range("F9").Value = "ORBITA 12 ROTOR WATCH `nWINDING CABINET WOOD `nOrbita"
Set-BlankCell "D9"
Set-BlankCell "E9"
Set-BlankCell "L9"
Set-BlankCell "M9"
Set-BlankCell "O9"
Set-BlankCell "P9"
Set-BlankCell "Q9"
Set-BlankCell "R9"
Set-BlankCell "S9"
Set-BlankCell "T9"

# --- Row 10: JAEGER-LECOULTRE --------------------------------------------
$ws.Range("C10").Value = "JAEGER-LECOULTRE DESK CLOCK BRASS `nJaeger-LeCoultre"
$ws.Range("F10").Value = "JAEGER-LECOULTRE DESK CLOCK BRASS `nJaeger-LeCoultre"
Set-BlankCell "D10"
Set-BlankCell "E10"
Set-BlankCell "L10"
Set-BlankCell "M10"
Set-BlankCell "O10"
Set-BlankCell "P10"
Set-BlankCell "Q10"
Set-BlankCell "R10"
Set-BlankCell "S10"
Set-BlankCell "T10"
